$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new section-header row above the old row 7 (2018.83 ...),
#    pushing the 2018.* rows down by one (rows 7-10 -> 8-11).
$ws.Rows("7:7").Insert()

# 2) Fill the new row 7 with just a label in column A.
$ws.Range("A7").Value2 = "2019b_moed_b"

# 3) Update the header in H1: "קומבינטוריקה כללי" -> "קומבינטוריקה מחרוזות ונסיגה"
$ws.Range("H1").Value2 = "קומבינטוריקה מחרוזות ונסיגה"

# 4) Widen column A to fit the new label.
$ws.Columns("A:A").ColumnWidth = 14.43

# 5) Update the active selection to H2 (matches the saved view state).
[void]$ws.Range("H2").Select()
